$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Set the Status column value for row 12 ("Pronto" = "Ready/Done")
$ws.Range("E12").Value = "Pronto"

# Update the active selection on the sheet from E14 to E13
$ws.Range("E13").Select()
